$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 189
$ws.Cells.Item(189, 2).Value = 7952779
$ws.Cells.Item(189, 5).Value = "Zrinjski Mostar"
$ws.Cells.Item(189, 6).Value = "FK Tuzla City"
$ws.Cells.Item(189, 7).Value = 4
$ws.Cells.Item(189, 8).Value = 0
$ws.Cells.Item(189, 9).Value = 2
$ws.Cells.Item(189, 10).Value = 0
$ws.Cells.Item(189, 11).Value = "H"
$ws.Cells.Item(189, 12).Value = 1.25
$ws.Cells.Item(189, 13).Value = 5.75
$ws.Cells.Item(189, 14).Value = 7
$ws.Cells.Item(189, 15).Value = 1.055
$ws.Cells.Item(189, 16).Value = 13
$ws.Cells.Item(189, 17).Value = 17
$ws.Cells.Item(189, 18).Value = -3.5
$ws.Cells.Item(189, 19).Value = 1.975
$ws.Cells.Item(189, 20).Value = 1.825
$ws.Cells.Item(189, 21).Value = 4.75
$ws.Cells.Item(189, 22).Value = 1.825
$ws.Cells.Item(189, 23).Value = 1.975
$ws.Cells.Item(189, 24).Value = 0.05499999999999994
$ws.Cells.Item(189, 25).Value = -1
$ws.Cells.Item(189, 26).Value = -1
$ws.Cells.Item(189, 27).Value = 0.9750000000000001
$ws.Cells.Item(189, 28).Value = -1
$ws.Cells.Item(189, 29).Value = -1
$ws.Cells.Item(189, 30).Value = 0.9750000000000001

# Row 191
$ws.Cells.Item(191, 2).Value = 7952780
$ws.Cells.Item(191, 5).Value = "Velez Mostar"
$ws.Cells.Item(191, 6).Value = "GOSK Gabela"
$ws.Cells.Item(191, 7).Value = 3
$ws.Cells.Item(191, 8).Value = 3
$ws.Cells.Item(191, 9).Value = 1
$ws.Cells.Item(191, 10).Value = 1
$ws.Cells.Item(191, 11).Value = "D"
$ws.Cells.Item(191, 12).Value = 1.4
$ws.Cells.Item(191, 13).Value = 4
$ws.Cells.Item(191, 14).Value = 7
$ws.Cells.Item(191, 15).Value = 1.363
$ws.Cells.Item(191, 16).Value = 4.2
$ws.Cells.Item(191, 17).Value = 8
$ws.Cells.Item(191, 18).Value = -1.5
$ws.Cells.Item(191, 19).Value = 2
$ws.Cells.Item(191, 20).Value = 1.8
$ws.Cells.Item(191, 21).Value = 2.75
$ws.Cells.Item(191, 22).Value = 1.825
$ws.Cells.Item(191, 23).Value = 1.975
$ws.Cells.Item(191, 24).Value = -1
$ws.Cells.Item(191, 25).Value = 3.2
$ws.Cells.Item(191, 26).Value = -1
$ws.Cells.Item(191, 27).Value = -1
$ws.Cells.Item(191, 28).Value = 0.8
$ws.Cells.Item(191, 29).Value = 0.825
$ws.Cells.Item(191, 30).Value = -1

# Row 192
$ws.Cells.Item(192, 2).Value = 7952778
$ws.Cells.Item(192, 5).Value = "Sloga"
$ws.Cells.Item(192, 6).Value = "Siroki Brijeg"
$ws.Cells.Item(192, 7).Value = 2
$ws.Cells.Item(192, 8).Value = 3
$ws.Cells.Item(192, 9).Value = 2
$ws.Cells.Item(192, 10).Value = 2
$ws.Cells.Item(192, 11).Value = "A"
$ws.Cells.Item(192, 12).Value = 1.727
$ws.Cells.Item(192, 13).Value = 3.75
$ws.Cells.Item(192, 14).Value = 3.75
$ws.Cells.Item(192, 15).Value = 1.7
$ws.Cells.Item(192, 16).Value = 3.9
$ws.Cells.Item(192, 17).Value = 3.9
$ws.Cells.Item(192, 18).Value = -0.75
$ws.Cells.Item(192, 19).Value = 1.975
$ws.Cells.Item(192, 20).Value = 1.825
$ws.Cells.Item(192, 21).Value = 2.25
$ws.Cells.Item(192, 22).Value = 1.8
$ws.Cells.Item(192, 23).Value = 2
$ws.Cells.Item(192, 24).Value = -1
$ws.Cells.Item(192, 25).Value = -1
$ws.Cells.Item(192, 26).Value = 2.9
$ws.Cells.Item(192, 27).Value = -1
$ws.Cells.Item(192, 28).Value = 0.825
$ws.Cells.Item(192, 29).Value = 0.8
$ws.Cells.Item(192, 30).Value = -1

# Row 193
$ws.Cells.Item(193, 2).Value = 7952776
$ws.Cells.Item(193, 5).Value = "FK Sarajevo"
$ws.Cells.Item(193, 6).Value = "NK Posusje"
$ws.Cells.Item(193, 7).Value = 1
$ws.Cells.Item(193, 8).Value = 1
$ws.Cells.Item(193, 9).Value = 0
$ws.Cells.Item(193, 10).Value = 0
$ws.Cells.Item(193, 11).Value = "D"
$ws.Cells.Item(193, 12).Value = 1.571
$ws.Cells.Item(193, 13).Value = 3.4
$ws.Cells.Item(193, 14).Value = 5.5
$ws.Cells.Item(193, 15).Value = 1.363
$ws.Cells.Item(193, 16).Value = 3.9
$ws.Cells.Item(193, 17).Value = 8
$ws.Cells.Item(193, 18).Value = -1.25
$ws.Cells.Item(193, 19).Value = 1.85
$ws.Cells.Item(193, 20).Value = 1.95
$ws.Cells.Item(193, 21).Value = 2.75
$ws.Cells.Item(193, 22).Value = 1.925
$ws.Cells.Item(193, 23).Value = 1.875
$ws.Cells.Item(193, 24).Value = -1
$ws.Cells.Item(193, 25).Value = 2.9
$ws.Cells.Item(193, 26).Value = -1
$ws.Cells.Item(193, 27).Value = -1
$ws.Cells.Item(193, 28).Value = 0.95
$ws.Cells.Item(193, 29).Value = -1
$ws.Cells.Item(193, 30).Value = 0.875

# Row 194
$ws.Cells.Item(194, 2).Value = 7952781
$ws.Cells.Item(194, 5).Value = "Zvijezda 09"
$ws.Cells.Item(194, 6).Value = "Zeljeznicar"
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 5
$ws.Cells.Item(194, 9).Value = 0
$ws.Cells.Item(194, 10).Value = 1
$ws.Cells.Item(194, 11).Value = "A"
$ws.Cells.Item(194, 12).Value = 2.15
$ws.Cells.Item(194, 13).Value = 3.25
$ws.Cells.Item(194, 14).Value = 2.9
$ws.Cells.Item(194, 15).Value = 3.6
$ws.Cells.Item(194, 16).Value = 3.4
$ws.Cells.Item(194, 17).Value = 1.85
$ws.Cells.Item(194, 18).Value = 0.5
$ws.Cells.Item(194, 19).Value = 1.875
$ws.Cells.Item(194, 20).Value = 1.925
$ws.Cells.Item(194, 21).Value = 2.5
$ws.Cells.Item(194, 22).Value = 1.975
$ws.Cells.Item(194, 23).Value = 1.825
$ws.Cells.Item(194, 24).Value = -1
$ws.Cells.Item(194, 25).Value = -1
$ws.Cells.Item(194, 26).Value = 0.8500000000000001
$ws.Cells.Item(194, 27).Value = -1
$ws.Cells.Item(194, 28).Value = 0.925
$ws.Cells.Item(194, 29).Value = 0.9750000000000001
$ws.Cells.Item(194, 30).Value = -1
